$d = $word.ActiveDocument

# 1. Replace the leading title text "Harvard DataFest2018" with "GIS Institute 2023".
#    (Keeps it within the same run/formatting; the trailing ": Geospatial Data
#    Management with " text is left untouched.)
$find = $d.Content.Find
$find.Execute("Harvard DataFest2018", $true, $false, $false, $false, $false, `
               $true, 1, $false, "GIS Institute 2023", 2)

# 2. Remove the old (hidden) "_GoBack" bookmark wherever it currently sits
#    (it was at the very end of the document, right after the last hyperlink
#    text in this file).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# 3. Re-insert the "_GoBack" bookmark immediately after the newly inserted
#    "GIS Institute 2023" text (collapsed/zero-length, matching Word's usual
#    "last edit position" bookmark).
$titleRange = $d.Content
$titleRange.Find.Execute("GIS Institute 2023")
$titleRange.Collapse(0)
$d.Bookmarks.Add("_GoBack", $titleRange)
